# Update the "timestamp" column (O) for all data rows (2-73) from the
# old scrape time "2022-07-19 07:02:16" to the new scrape time
# "2022-07-19 20:57:47".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-07-19 20:57:47"

for ($row = 2; $row -le 73; $row++) {
    $ws.Cells.Item($row, 15).Value = $newTimestamp
}
